$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 74 (date serial 45630, 2024-12-04).
# Append two more rows continuing the same daily series, with the same
# metric values as the last existing row, for 2024-12-05 and 2024-12-06
# (date serials 45631 and 45632).

# Row 75: duplicate row 74, then update the date in column A.
$ws.Range("A74:J74").Copy($ws.Range("A75:J75"))
$ws.Range("A75").Value2 = 45631

# Row 76: duplicate row 74, then update the date in column A.
$ws.Range("A74:J74").Copy($ws.Range("A76:J76"))
$ws.Range("A76").Value2 = 45632
